$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2159090909090909
$ws.Range("C2").Value = 0.5303030303030303
$ws.Range("J2").Value = 0.01136363636363636
$ws.Range("P2").Value = 0.1287878787878788
$ws.Range("S2").Value = 0.1136363636363636
$ws.Range("B3").Value = 0.006896551724137931
$ws.Range("C3").Value = 0.02758620689655172
$ws.Range("J3").Value = 0.02758620689655172
$ws.Range("P3").Value = 0.7310344827586207
$ws.Range("S3").Value = 0.2068965517241379
$ws.Range("P4").Value = 0.7575757575757576
$ws.Range("S4").Value = 0.2424242424242424
$ws.Range("B6").Value = 0.06779661016949153
$ws.Range("D6").Value = 0.01271186440677966
$ws.Range("F6").Value = 0.08050847457627118
$ws.Range("J6").Value = 0.2584745762711864
$ws.Range("O6").Value = 0.0211864406779661
$ws.Range("Q6").Value = 0.1271186440677966
$ws.Range("R6").Value = 0.05084745762711865
$ws.Range("S6").Value = 0.3813559322033898
$ws.Range("B7").Value = 0.1015228426395939
$ws.Range("D7").Value = 0.02538071065989848
$ws.Range("F7").Value = 0.07106598984771574
$ws.Range("J7").Value = 0.1472081218274112
$ws.Range("O7").Value = 0.01522842639593909
$ws.Range("Q7").Value = 0.116751269035533
$ws.Range("R7").Value = 0.06091370558375635
$ws.Range("S7").Value = 0.4619289340101523
$ws.Range("B8").Value = 0.1180904522613065
$ws.Range("D8").Value = 0.02010050251256281
$ws.Range("E8").Value = 0.002512562814070352
$ws.Range("F8").Value = 0.08040201005025126
$ws.Range("J8").Value = 0.1055276381909548
$ws.Range("O8").Value = 0.02010050251256281
$ws.Range("Q8").Value = 0.1582914572864322
$ws.Range("R8").Value = 0.08793969849246232
$ws.Range("S8").Value = 0.407035175879397
$ws.Range("B9").Value = 0.06369426751592357
$ws.Range("D9").Value = 0.01273885350318471
$ws.Range("F9").Value = 0.07643312101910828
$ws.Range("J9").Value = 0.1019108280254777
$ws.Range("O9").Value = 0.01910828025477707
$ws.Range("Q9").Value = 0.2165605095541401
$ws.Range("R9").Value = 0.05095541401273886
$ws.Range("S9").Value = 0.4585987261146497
$ws.Range("B10").Value = 0.1041069723018147
$ws.Range("D10").Value = 0.01432664756446991
$ws.Range("F10").Value = 0.0830945558739255
$ws.Range("J10").Value = 0.113658070678128
$ws.Range("O10").Value = 0.01528175740210124
$ws.Range("Q10").Value = 0.1900668576886342
$ws.Range("R10").Value = 0.09264565425023878
$ws.Range("S10").Value = 0.3868194842406877
$ws.Range("G11").Value = 0.1395348837209302
$ws.Range("J11").Value = 0.09634551495016612
$ws.Range("K11").Value = 0.186046511627907
$ws.Range("L11").Value = 0.5614617940199336
$ws.Range("S11").Value = 0.01661129568106312
$ws.Range("G12").Value = 0.7816091954022989
$ws.Range("J12").Value = 0.1264367816091954
$ws.Range("K12").Value = 0.01724137931034483
$ws.Range("L12").Value = 0.02873563218390805
$ws.Range("S12").Value = 0.04597701149425287
$ws.Range("G13").Value = 0.6153846153846154
$ws.Range("J13").Value = 0.2307692307692308
$ws.Range("S13").Value = 0.1538461538461539
$ws.Range("F15").Value = 0.009259259259259259
$ws.Range("H15").Value = 0.1712962962962963
$ws.Range("I15").Value = 0.06944444444444445
$ws.Range("J15").Value = 0.3842592592592592
$ws.Range("K15").Value = 0.08333333333333333
$ws.Range("O15").Value = 0.07407407407407407
$ws.Range("S15").Value = 0.2083333333333333
$ws.Range("F16").Value = 0.0440251572327044
$ws.Range("H16").Value = 0.1949685534591195
$ws.Range("I16").Value = 0.06918238993710692
$ws.Range("J16").Value = 0.3459119496855346
$ws.Range("K16").Value = 0.1257861635220126
$ws.Range("M16").Value = 0.01886792452830189
$ws.Range("O16").Value = 0.05660377358490566
$ws.Range("S16").Value = 0.1446540880503145
$ws.Range("F17").Value = 0.02285714285714286
$ws.Range("H17").Value = 0.1771428571428571
$ws.Range("I17").Value = 0.1114285714285714
$ws.Range("J17").Value = 0.3771428571428572
$ws.Range("K17").Value = 0.1085714285714286
$ws.Range("M17").Value = 0.02285714285714286
$ws.Range("O17").Value = 0.06857142857142857
$ws.Range("S17").Value = 0.1114285714285714
$ws.Range("F18").Value = 0.04242424242424243
$ws.Range("H18").Value = 0.1636363636363636
$ws.Range("I18").Value = 0.1090909090909091
$ws.Range("J18").Value = 0.4
$ws.Range("K18").Value = 0.07878787878787878
$ws.Range("M18").Value = 0.02424242424242424
$ws.Range("O18").Value = 0.103030303030303
$ws.Range("S18").Value = 0.07878787878787878
$ws.Range("F19").Value = 0.02028218694885361
$ws.Range("H19").Value = 0.2142857142857143
$ws.Range("I19").Value = 0.06701940035273368
$ws.Range("J19").Value = 0.3492063492063492
$ws.Range("K19").Value = 0.1322751322751323
$ws.Range("M19").Value = 0.02204585537918871
$ws.Range("N19").Value = 0.0008818342151675485
$ws.Range("O19").Value = 0.08201058201058201
$ws.Range("S19").Value = 0.1119929453262787
